$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 687
    $ws.Range("F3").Value = 528
    $ws.Range("F8").Value = 3355
    $ws.Range("F10").Value = 124
}
